$d = $word.ActiveDocument

function Add-Split($pos, $name) {
    # Inserting (and immediately deleting) a bookmark at a text position forces
    # the run that spans that position to be split into two runs at that exact
    # offset, without merging/altering any of the surrounding runs.
    # (NOTE: $d is captured from the enclosing scope rather than passed as a
    # parameter -- passing a COM object as a positional arg together with a
    # parenthesized expression confuses this shell's argument binder.)
    $r = $d.Range($pos, $pos)
    $d.Bookmarks.Add($name, $r)
    $d.Bookmarks.Item($name).Delete()
}

# ---------------------------------------------------------------------------
# Change 1: "Generate a PDB SDP Map from SDP Scores in MySQL:" heading
#   -> "Generate a PDB Map from in MySQL:" (with the _GoBack bookmark moved to
#      sit right after "in")
# ---------------------------------------------------------------------------

$paras = $d.Paragraphs
$headingIndex = -1
for ($i = 1; $i -le $paras.Count; $i++) {
    if ($paras.Item($i).Range.Text -like "Generate a PDB SDP Map*") {
        $headingIndex = $i
        break
    }
}

$p1 = $paras.Item($headingIndex)
$base1 = $p1.Range.Start

# Original text: "Generate a PDB SDP Map from SDP Scores in MySQL:"
# Delete "SDP Scores " (rightmost first, so earlier offsets stay valid)
$d.Range($base1+28, $base1+39).Delete()
# Delete "SDP " (right after "PDB ")
$d.Range($base1+15, $base1+19).Delete()

# Text is now: "Generate a PDB Map from in MySQL:"
# Split into the desired run boundaries:
#   "Generate a PDB Map " | "from " | "in" | [[bookmark]] | " MySQL:"
Add-Split ($base1+19) "TempSplit1"
Add-Split ($base1+24) "TempSplit2"

# Move the _GoBack bookmark to right after "in" (offset 26)
$bmPos = $d.Range($base1+26, $base1+26)
$d.Bookmarks.Add("_GoBack", $bmPos)

# ---------------------------------------------------------------------------
# Change 2: "movie.roll 1,180,1,axis=y" + ";" (two runs) -> a single run
#   "movie.roll 1,180,1,axis=y;"
# ---------------------------------------------------------------------------

$paras = $d.Paragraphs
$movieIndex = -1
for ($i = 1; $i -le $paras.Count; $i++) {
    if ($paras.Item($i).Range.Text -like "orient;mset*") {
        $movieIndex = $i
        break
    }
}

$p2 = $paras.Item($movieIndex)
$base2 = $p2.Range.Start

# Original text: "orient;mset 1 x360; movie.roll 1,180,1,axis=y; movie.roll 181,360,1,axis=x;"
# Replace "movie.roll 1,180,1,axis=y;" with itself - this merges the previously
# separate "movie.roll 1,180,1,axis=y" and ";" runs into one run.
$narrow = $d.Range($base2+20, $base2+46)
$narrow.Find.Execute("movie.roll 1,180,1,axis=y;", $false, $false, $false, $false, $false, $true, 1, $false, "movie.roll 1,180,1,axis=y;", 2)

# Re-establish the other original run boundaries (space / text / space / text),
# splitting from right to left so each split's whitespace handling resolves
# correctly.
Add-Split ($base2+47) "TempSplitR"
Add-Split ($base2+46) "TempSplitM"
Add-Split ($base2+20) "TempSplitL"
